$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the mismatched Name/City values in row 5:
# B5 currently holds "Mbanga" (a city), C5 currently holds "Lapiro" (a name).
# They were entered in the wrong columns, so swap them.
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$ws.Range("B5").Value2 = $c5
$ws.Range("C5").Value2 = $b5

# Update the selected cell / active cell shown in the saved view.
$ws.Range("C5").Select()
